$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B, shifting C:H left into B:G
$ws.Range("B:B").Delete()

# Update header row text (append ".jamais.jamais")
$ws.Range("B1").Value = "Stable accommodation.jamais.jamais"
$ws.Range("C1").Value = "Unstable accommodation and/or homeless.jamais.jamais"
$ws.Range("D1").Value = "In detention.jamais.jamais"
$ws.Range("E1").Value = "Other.jamais.jamais"
$ws.Range("F1").Value = "Not known / missing.jamais.jamais"
$ws.Range("G1").Value = "Total.jamais.jamais"
